# Lin EIS Sens Study - Warburg Params
# Adds two new parameter rows ("mD" and "tauF") to the Parameters sheet,
# just above the "Normalized reaction rate coefficient(s)" row (original
# row 58), and hides the Sectioning column/header rows at the top of the
# sheet. Also nudges sheet view/selection to match the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# ---------------------------------------------------------------------
# 1) Insert the "mD" row at row 58 (everything from old row 58 on shifts
#    down by one). Borrow the formatting from row 57 ("Dsref"), which
#    already carries the B/C/D/E/F/G style pattern (6,6,6,37,22,8) that
#    the new row needs.
# ---------------------------------------------------------------------
$ws.Rows.Item(58).Insert()
$ws.Range("B57:G57").Copy()
$ws.Range("B58:G58").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Insert the "tauF" row at row 60 (pushes the row that is currently at
#    60 - originally row 59 - down to row 61). Borrow formatting from row
#    59 (originally row 58, "nF"), which already has the desired style
#    pattern (6,6,6,22,22,8).
# ---------------------------------------------------------------------
$ws.Rows.Item(60).Insert()
$ws.Range("B59:G59").Copy()
$ws.Range("B60:G60").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Populate the new "tauF" row (row 60) first so the new shared
#    strings land in the same order as the authored workbook.
# ---------------------------------------------------------------------
$ws.Range("B60").Value = "Solid-phase diffusivity CPE-integrator time const."
$ws.Range("C60").Value = "tauF"
$ws.Range("D60").Value = "\tau_\mathrm{f}"
$ws.Range("E60").Value = 100
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = "s"

# ---------------------------------------------------------------------
# 4) Populate the new "mD" row (row 58).
# ---------------------------------------------------------------------
$ws.Range("B58").Value = "Empirical multiplicity of solid diffusivity"
$ws.Range("C58").Value = "mD"
$ws.Range("D58").Value = "m_\mathrm{D}"
$ws.Range("E58").Value = 1
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = "unitless"

# ---------------------------------------------------------------------
# 5) Hide the "Sectioning" column (A) and the top title/header rows
#    (1 and 2), matching the authored edit.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).Hidden = $true
$ws.Rows.Item(1).Hidden = $true
$ws.Rows.Item(2).Hidden = $true

# ---------------------------------------------------------------------
# 6) Update the view: scroll so row 16 is the top-most visible row,
#    zoom to 110%, and select E27 (matches authored sheetView/selection).
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 110
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("E27").Select()

Write-Output "done"
